$wb = $excel.ActiveWorkbook

$targets = @{
    "2025" = @{ "B2" = 779.1630278080007; "E2" = 23185.90077278444; "I2" = 12940.225083824; "L2" = 38819.6236028304; "M2" = 8707.039609781334; "N2" = 5744.096799946036; "O2" = 5589.635007435185 }
    "2030" = @{ "B2" = 4570.048142274309; "E2" = 44884.8447102553; "I2" = 35374.31877773289; "L2" = 53573.26199887139; "M2" = 18074.35048528854; "N2" = 8506.387057313099; "O2" = 9664.202924111112 }
    "2035" = @{ "A2" = 2289.569120990697; "B2" = 6421.511730469836; "E2" = 53838.19196405857; "I2" = 47405.14060738889; "L2" = 53573.26199887139; "M2" = 20935.70896779413; "N2" = 12148.8294540729; "O2" = 11823.7365588148 }
    "2040" = @{ "A2" = 2289.569120990697; "B2" = 6421.511730469836; "E2" = 53838.19196405857; "I2" = 47405.14060738889; "L2" = 53573.26199887139; "M2" = 20935.70896779413; "N2" = 12233.16180906476; "O2" = 11823.7365588148 }
    "2045" = @{ "A2" = 5042.298495067215; "B2" = 6421.511730469836; "E2" = 53838.19196405857; "I2" = 47405.14060738889; "L2" = 53573.26199887139; "M2" = 20935.70896779413; "N2" = 12666.29850282887; "O2" = 13702.58732167914 }
    "2050" = @{ "A2" = 5042.298495067215; "B2" = 6421.511730469836; "E2" = 53838.19196405857; "I2" = 47405.14060738889; "L2" = 53573.26199887139; "M2" = 20935.70896779413; "N2" = 12666.29850282887; "O2" = 13702.58732167914 }
}

foreach ($sheetName in $targets.Keys) {
    $sName = [string]$sheetName
    $ws = $wb.Worksheets.Item($sName)
    $cellVals = $targets[$sheetName]
    foreach ($cellRef in $cellVals.Keys) {
        $ws.Range($cellRef).Value = $cellVals[$cellRef]
    }
}
